$d = $word.ActiveDocument

$d.Content.Find.Execute("41 730 000,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "41 770 000,00", 2)

$d.Content.Find.Execute("3 755 700,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3 759 300,00", 2)
